# Update imputed values in column A and C to reflect the re-run of the
# RandomForest algorithm (commit: "Update Name of Algo").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "A3"   = -22.139
    "A14"  = -21.5877
    "C15"  = -14.16499999999999
    "A16"  = -22.10020000000001
    "A21"  = -20.19219999999999
    "C21"  = -12.1812
    "C22"  = -11.6268
    "A23"  = -19.95949999999998
    "C24"  = -13.47399999999999
    "A25"  = -21.57519999999998
    "A26"  = -21.18959999999997
    "C27"  = -12.2965
    "C28"  = -13.53209999999999
    "A29"  = -20.75329999999997
    "C36"  = -13.3351
    "C39"  = -13.4066
    "A40"  = -19.9843
    "C45"  = -13.96439999999999
    "C48"  = -11.9996
    "C49"  = -13.24290000000001
    "C52"  = -10.6379
    "A53"  = -21.98439999999999
    "C53"  = -12.68299999999999
    "C54"  = -13.4969
    "A57"  = -22.17140000000001
    "C57"  = -14.3953
    "A59"  = -22.42610000000001
    "A65"  = -21.87619999999999
    "A69"  = -21.6085
    "C70"  = -11.9309
    "C71"  = -11.1282
    "A79"  = -20.634
    "A83"  = -21.9176
    "C86"  = -13.89139999999999
    "C87"  = -12.99520000000001
    "C89"  = -13.5628
    "A91"  = -20.46169999999998
    "A93"  = -21.25849999999999
    "A100" = -22.20870000000001
    "C101" = -13.48639999999999
    "A103" = -21.89179999999999
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
